$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"3"
$ws.Range("F2").Value = [double]"1"
$ws.Range("G2").Value = [double]"5.530908"
$ws.Range("H2").Value = [double]"16.592724"
$ws.Range("I2").Value = [double]"0.06269514438603573"
$ws.Range("J2").Value = [double]"0.06269514438603574"
$ws.Range("K2").Value = [double]"3"
$ws.Range("L2").Value = [double]"1"
$ws.Range("M2").Value = [double]"7.107333666666666"
$ws.Range("N2").Value = [double]"21.322001"
$ws.Range("O2").Value = [double]"0.7373665550576455"
$ws.Range("P2").Value = [double]"0.7373665550576454"
$ws.Range("Q2").Value = [double]"39.310008635636"
$ws.Range("R2").Value = [double]"353.790077720724"
$ws.Range("S2").Value = [double]"0.04622930263477285"
$ws.Range("T2").Value = [double]"0.04622930263477285"
$ws.Range("E3").Value = [double]"3"
$ws.Range("F3").Value = [double]"1"
$ws.Range("G3").Value = [double]"5.530908"
$ws.Range("H3").Value = [double]"16.592724"
$ws.Range("I3").Value = [double]"0.06269514438603573"
$ws.Range("J3").Value = [double]"0.06269514438603574"
$ws.Range("O3").Value = [double]"0.1688878844614928"
$ws.Range("P3").Value = [double]"0.1688878844614928"
$ws.Range("Q3").Value = [double]"9.003641609588"
$ws.Range("R3").Value = [double]"81.032774486292"
$ws.Range("S3").Value = [double]"0.01058845030136541"
$ws.Range("T3").Value = [double]"0.01058845030136541"
$ws.Range("E4").Value = [double]"3"
$ws.Range("F4").Value = [double]"1"
$ws.Range("G4").Value = [double]"5.530908"
$ws.Range("H4").Value = [double]"16.592724"
$ws.Range("I4").Value = [double]"0.06269514438603573"
$ws.Range("J4").Value = [double]"0.06269514438603574"
$ws.Range("M4").Value = [double]"0.8135026666666666"
$ws.Range("N4").Value = [double]"2.440508"
$ws.Range("O4").Value = [double]"0.08439869112428164"
$ws.Range("P4").Value = [double]"0.08439869112428162"
$ws.Range("Q4").Value = [double]"4.499408407088"
$ws.Range("R4").Value = [double]"40.494675663792"
$ws.Range("S4").Value = [double]"0.005291388126029269"
$ws.Range("T4").Value = [double]"0.00529138812602927"
$ws.Range("E5").Value = [double]"3"
$ws.Range("F5").Value = [double]"1"
$ws.Range("G5").Value = [double]"5.530908"
$ws.Range("H5").Value = [double]"16.592724"
$ws.Range("I5").Value = [double]"0.06269514438603573"
$ws.Range("J5").Value = [double]"0.06269514438603574"
$ws.Range("M5").Value = [double]"0.09009266666666667"
$ws.Range("N5").Value = [double]"0.270278"
$ws.Range("O5").Value = [double]"0.009346869356580103"
$ws.Range("P5").Value = [double]"0.009346869356580103"
$ws.Range("Q5").Value = [double]"0.498294250808"
$ws.Range("R5").Value = [double]"4.484648257272"
$ws.Range("S5").Value = [double]"0.0005860033238682024"
$ws.Range("T5").Value = [double]"0.0005860033238682025"
$ws.Range("I6").Value = [double]"0.6763730500901992"
$ws.Range("J6").Value = [double]"0.6763730500901993"
$ws.Range("K6").Value = [double]"3"
$ws.Range("L6").Value = [double]"1"
$ws.Range("M6").Value = [double]"7.107333666666666"
$ws.Range("N6").Value = [double]"21.322001"
$ws.Range("O6").Value = [double]"0.7373665550576455"
$ws.Range("P6").Value = [double]"0.7373665550576454"
$ws.Range("Q6").Value = [double]"424.0875541532251"
$ws.Range("R6").Value = [double]"3816.787987379026"
$ws.Range("S6").Value = [double]"0.4987348658788425"
$ws.Range("T6").Value = [double]"0.4987348658788425"
$ws.Range("I7").Value = [double]"0.6763730500901992"
$ws.Range("J7").Value = [double]"0.6763730500901993"
$ws.Range("O7").Value = [double]"0.1688878844614928"
$ws.Range("P7").Value = [double]"0.1688878844614928"
$ws.Range("S7").Value = [double]"0.1142312135365011"
$ws.Range("T7").Value = [double]"0.1142312135365011"
$ws.Range("I8").Value = [double]"0.6763730500901992"
$ws.Range("J8").Value = [double]"0.6763730500901993"
$ws.Range("M8").Value = [double]"0.8135026666666666"
$ws.Range("N8").Value = [double]"2.440508"
$ws.Range("O8").Value = [double]"0.08439869112428164"
$ws.Range("P8").Value = [double]"0.08439869112428162"
$ws.Range("Q8").Value = [double]"48.54089766768977"
$ws.Range("R8").Value = [double]"436.868079009208"
$ws.Range("S8").Value = [double]"0.05708500013935099"
$ws.Range("T8").Value = [double]"0.05708500013935099"
$ws.Range("I9").Value = [double]"0.6763730500901992"
$ws.Range("J9").Value = [double]"0.6763730500901993"
$ws.Range("K9").Value = [double]"2"
$ws.Range("L9").Value = [double]"0.6666666666666666"
$ws.Range("M9").Value = [double]"0.09009266666666667"
$ws.Range("N9").Value = [double]"0.270278"
$ws.Range("O9").Value = [double]"0.009346869356580103"
$ws.Range("P9").Value = [double]"0.009346869356580103"
$ws.Range("Q9").Value = [double]"5.375740108136444"
$ws.Range("R9").Value = [double]"48.381660973228"
$ws.Range("S9").Value = [double]"0.006321970535504702"
$ws.Range("T9").Value = [double]"0.006321970535504703"
$ws.Range("G10").Value = [double]"23.01602366666667"
$ws.Range("H10").Value = [double]"69.04807100000001"
$ws.Range("I10").Value = [double]"0.2608962085382874"
$ws.Range("J10").Value = [double]"0.2608962085382874"
$ws.Range("K10").Value = [double]"3"
$ws.Range("L10").Value = [double]"1"
$ws.Range("M10").Value = [double]"7.107333666666666"
$ws.Range("N10").Value = [double]"21.322001"
$ws.Range("O10").Value = [double]"0.7373665550576455"
$ws.Range("P10").Value = [double]"0.7373665550576454"
$ws.Range("Q10").Value = [double]"163.5825598788968"
$ws.Range("R10").Value = [double]"1472.243038910071"
$ws.Range("S10").Value = [double]"0.1923761385174781"
$ws.Range("T10").Value = [double]"0.192376138517478"
$ws.Range("G11").Value = [double]"23.01602366666667"
$ws.Range("H11").Value = [double]"69.04807100000001"
$ws.Range("I11").Value = [double]"0.2608962085382874"
$ws.Range("J11").Value = [double]"0.2608962085382874"
$ws.Range("O11").Value = [double]"0.1688878844614928"
$ws.Range("P11").Value = [double]"0.1688878844614928"
$ws.Range("Q11").Value = [double]"37.46727090243811"
$ws.Range("R11").Value = [double]"337.205438121943"
$ws.Range("S11").Value = [double]"0.04406220872405583"
$ws.Range("T11").Value = [double]"0.04406220872405581"
$ws.Range("G12").Value = [double]"23.01602366666667"
$ws.Range("H12").Value = [double]"69.04807100000001"
$ws.Range("I12").Value = [double]"0.2608962085382874"
$ws.Range("J12").Value = [double]"0.2608962085382874"
$ws.Range("M12").Value = [double]"0.8135026666666666"
$ws.Range("N12").Value = [double]"2.440508"
$ws.Range("O12").Value = [double]"0.08439869112428164"
$ws.Range("P12").Value = [double]"0.08439869112428162"
$ws.Range("Q12").Value = [double]"18.72359662889644"
$ws.Range("R12").Value = [double]"168.512369660068"
$ws.Range("S12").Value = [double]"0.02201929851991909"
$ws.Range("T12").Value = [double]"0.02201929851991909"
$ws.Range("G13").Value = [double]"23.01602366666667"
$ws.Range("H13").Value = [double]"69.04807100000001"
$ws.Range("I13").Value = [double]"0.2608962085382874"
$ws.Range("J13").Value = [double]"0.2608962085382874"
$ws.Range("K13").Value = [double]"2"
$ws.Range("L13").Value = [double]"0.6666666666666666"
$ws.Range("M13").Value = [double]"0.09009266666666667"
$ws.Range("N13").Value = [double]"0.270278"
$ws.Range("O13").Value = [double]"0.009346869356580103"
$ws.Range("P13").Value = [double]"0.009346869356580103"
$ws.Range("Q13").Value = [double]"2.073574948193111"
$ws.Range("R13").Value = [double]"18.662174533738"
$ws.Range("S13").Value = [double]"0.002438562776834451"
$ws.Range("T13").Value = [double]"0.002438562776834451"
$ws.Range("G14").Value = [double]"0.003140333333333334"
$ws.Range("H14").Value = [double]"0.009421000000000001"
$ws.Range("I14").Value = [double]"3.559698547754079E-05"
$ws.Range("J14").Value = [double]"3.55969854775408E-05"
$ws.Range("K14").Value = [double]"3"
$ws.Range("L14").Value = [double]"1"
$ws.Range("M14").Value = [double]"7.107333666666666"
$ws.Range("N14").Value = [double]"21.322001"
$ws.Range("O14").Value = [double]"0.7373665550576455"
$ws.Range("P14").Value = [double]"0.7373665550576454"
$ws.Range("Q14").Value = [double]"0.02231939682455556"
$ws.Range("R14").Value = [double]"0.200874571421"
$ws.Range("S14").Value = [double]"2.624802655201129E-05"
$ws.Range("T14").Value = [double]"2.624802655201129E-05"
$ws.Range("G15").Value = [double]"0.003140333333333334"
$ws.Range("H15").Value = [double]"0.009421000000000001"
$ws.Range("I15").Value = [double]"3.559698547754079E-05"
$ws.Range("J15").Value = [double]"3.55969854775408E-05"
$ws.Range("O15").Value = [double]"0.1688878844614928"
$ws.Range("P15").Value = [double]"0.1688878844614928"
$ws.Range("Q15").Value = [double]"0.005112078499222223"
$ws.Range("R15").Value = [double]"0.046008706493"
$ws.Range("S15").Value = [double]"6.011899570508348E-06"
$ws.Range("T15").Value = [double]"6.011899570508347E-06"
$ws.Range("G16").Value = [double]"0.003140333333333334"
$ws.Range("H16").Value = [double]"0.009421000000000001"
$ws.Range("I16").Value = [double]"3.559698547754079E-05"
$ws.Range("J16").Value = [double]"3.55969854775408E-05"
$ws.Range("M16").Value = [double]"0.8135026666666666"
$ws.Range("N16").Value = [double]"2.440508"
$ws.Range("O16").Value = [double]"0.08439869112428164"
$ws.Range("P16").Value = [double]"0.08439869112428162"
$ws.Range("Q16").Value = [double]"0.002554669540888889"
$ws.Range("R16").Value = [double]"0.022992025868"
$ws.Range("S16").Value = [double]"3.004338982274505E-06"
$ws.Range("T16").Value = [double]"3.004338982274505E-06"
$ws.Range("G17").Value = [double]"0.003140333333333334"
$ws.Range("H17").Value = [double]"0.009421000000000001"
$ws.Range("I17").Value = [double]"3.559698547754079E-05"
$ws.Range("J17").Value = [double]"3.55969854775408E-05"
$ws.Range("K17").Value = [double]"2"
$ws.Range("L17").Value = [double]"0.6666666666666666"
$ws.Range("M17").Value = [double]"0.09009266666666667"
$ws.Range("N17").Value = [double]"0.270278"
$ws.Range("O17").Value = [double]"0.009346869356580103"
$ws.Range("P17").Value = [double]"0.009346869356580103"
$ws.Range("Q17").Value = [double]"0.0002829210042222222"
$ws.Range("R17").Value = [double]"0.002546289038"
$ws.Range("S17").Value = [double]"3.32720372746653E-07"
$ws.Range("T17").Value = [double]"3.32720372746653E-07"
